$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "314.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.87%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.41"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.85%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.124"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.05%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08202"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.82%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.967"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.31%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.246"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.75%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9315"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.24%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.20%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1972"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.24%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09130"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.66%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03523"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.13%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09805"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.22%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001395"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.05%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005993"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.01%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.660"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.74%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.262"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.07%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.71%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3464"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.07%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1285"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.82%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.884"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.88%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2445"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.54%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04317"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.89%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001219"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.98%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004800"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "16.94%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001295"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.57%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003987"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.36%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02248"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.54%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05264"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.77%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007559"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.99%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009821"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.24%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1378"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.45%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002112"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.52%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009791"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.23%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006352"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.40%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000748"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.79%"

$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001196"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-25.54%"

$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002760"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.82%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002093"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.79%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001994"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.79%"
